$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: extend the thin-bottom-border strip into the new column L ---
$ws.Range("J3").Copy()
$ws.Range("L3").PasteSpecial(-4122)

# --- Row 4 (header): new year column 2021, same look as K4 (bold Times,
#     medium border) but only a bottom border, matching D4:J4's style ---
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 2021
$ws.Range("L4").Borders.Item(8).LineStyle = 0

# --- Row 5 (Small enterprises): new data point, non-bold, no border ---
$ws.Range("K4").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Value = 2.3
$ws.Range("L5").Font.Bold = $false
$ws.Range("L5").Borders.Item(8).LineStyle = 0
$ws.Range("L5").Borders.Item(9).LineStyle = 0

# --- Row 6 (Medium-sized enterprises): new data point, non-bold, bottom
#     border only (row 6 is the bottom of the table, thickBot) ---
$ws.Range("K4").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L6").Value = 1.3
$ws.Range("L6").Font.Bold = $false
$ws.Range("L6").Borders.Item(8).LineStyle = 0

# --- Move the active selection as recorded in the saved view state ---
[void]$ws.Range("O5").Select()
